$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to stay plain text (avoids Excel
# auto-converting numeric-looking strings like "1.004" into a Number), and
# without leaving a residual style on the cell.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = "29.192.03"
$ws.Range("E2").Value = "  +1.94%  "

# Row 3
$ws.Range("D3").Value = "1.910.72"
$ws.Range("E3").Value = "  +2.25%  "

# Row 4
Set-TextValue "D4" "1.004"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
Set-TextValue "D5" "327.98"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6
Set-TextValue "D6" "1.004"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7
Set-TextValue "D7" "0.4659"
$ws.Range("E7").Value = "  +0.57%  "

# Row 8
Set-TextValue "D8" "0.3930"
$ws.Range("E8").Value = "  +1.20%  "

# Row 9
Set-TextValue "D9" "46.97"
$ws.Range("E9").Value = "  +0.82%  "

# Row 10
$ws.Range("E10").Value = "  +1.36%  "

# Row 11
$ws.Range("E11").Value = "  +3.03%  "

# Row 12
Set-TextValue "D12" "22.38"
$ws.Range("E12").Value = "  +1.97%  "

# Row 13
$ws.Range("D13").Value = "1.975.19"
$ws.Range("E13").Value = "  +4.01%  "

# Row 14
Set-TextValue "D14" "7.152"
$ws.Range("E14").Value = "  +2.02%  "

# Row 15
Set-TextValue "D15" "5.798"
$ws.Range("E15").Value = "  +1.68%  "

# Row 16
Set-TextValue "D16" "0.06972"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
Set-TextValue "D17" "88.82"
$ws.Range("E17").Value = "  +0.88%  "

# Row 18
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
Set-TextValue "D19" "0.00001009"
$ws.Range("E19").Value = "  +0.64%  "

# Row 20
Set-TextValue "D20" "17.25"
$ws.Range("E20").Value = "  +2.56%  "

# Row 21
Set-TextValue "D21" "1.003"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").Value = "29.210.35"
$ws.Range("E22").Value = "  +1.96%  "

# Row 23
$ws.Range("E23").Value = "  +1.82%  "

# Row 24
Set-TextValue "D24" "11.09"
$ws.Range("E24").Value = "  +0.79%  "

# Row 25
$ws.Range("D25").Value = "2.186.68"
$ws.Range("E25").Value = "  +2.92%  "

# Row 26
Set-TextValue "D26" "2.062"
$ws.Range("E26").Value = "  -2.64%  "

# Row 27
Set-TextValue "D27" "156.39"
$ws.Range("E27").Value = "  +2.38%  "

# Row 28
$ws.Range("E28").Value = "  +1.79%  "

# Row 29
Set-TextValue "D29" "5.864"
$ws.Range("E29").Value = "  +1.43%  "

# Row 30
Set-TextValue "D30" "2.010"
$ws.Range("E30").Value = "  +1.24%  "

# Row 31
Set-TextValue "D31" "119.69"
$ws.Range("E31").Value = "  +0.31%  "

# Row 32
Set-TextValue "D32" "0.09419"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33
Set-TextValue "D33" "0.9255"
$ws.Range("E33").Value = "  +0.79%  "

# Row 34
Set-TextValue "D34" "5.371"
$ws.Range("E34").Value = "  +1.91%  "

# Row 35
Set-TextValue "D35" "1.347"
$ws.Range("E35").Value = "  +0.73%  "

# Row 36
$ws.Range("E36").Value = "  -1.96%  "

# Row 37
Set-TextValue "D37" "0.05864"
$ws.Range("E37").Value = "  +1.43%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.165"
$ws.Range("E38").Value = "  +1.14%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D39" "8.048"
$ws.Range("E39").Value = "  +3.67%  "

# Row 40
Set-TextValue "D40" "0.02100"
$ws.Range("E40").Value = "  -0.38%  "

# Row 41
Set-TextValue "D41" "0.5761"
$ws.Range("E41").Value = "  +2.36%  "

# Row 42
Set-TextValue "D42" "0.1814"
$ws.Range("E42").Value = "  +1.39%  "

# Row 43
Set-TextValue "D43" "10.04"
$ws.Range("E43").Value = "  +2.86%  "

# Row 44
Set-TextValue "D44" "12.10"
$ws.Range("E44").Value = "  +2.86%  "

# Row 45
Set-TextValue "D45" "0.5434"
$ws.Range("E45").Value = "  +2.26%  "

# Row 46
Set-TextValue "D46" "2.231"
$ws.Range("E46").Value = "  +5.72%  "

# Row 47
Set-TextValue "D47" "0.07103"
$ws.Range("E47").Value = "  -1.27%  "

# Row 48
Set-TextValue "D48" "1.884"
$ws.Range("E48").Value = "  +2.98%  "

# Row 49
Set-TextValue "D49" "2.571"
$ws.Range("E49").Value = "  +6.74%  "

# Row 50
Set-TextValue "D50" "112.56"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51
Set-TextValue "D51" "1.089"
$ws.Range("E51").Value = "  -4.77%  "
